# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" tracker sheet.
#
# Layout: col A = date (serial), B/C/D = Chase/Bryce/Zach win counts.
# The most-recent row is always formatted as a plain date ("YYYY-MM-DD"),
# while every earlier row uses the full timestamp format
# ("YYYY-MM-DD HH:MM:SS"). When a new row is appended, the previous last
# row reverts to the timestamp format and the new row gets the date-only
# format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 78
$newRow = $lastRow + 1

# Previous last row (78) goes back to the "regular" timestamp format.
$ws.Range("A" + $lastRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New last row (79) gets the "latest entry" date-only format plus today's data.
$ws.Range("A" + $newRow).NumberFormat = "YYYY-MM-DD"
$ws.Range("A" + $newRow).Value = 45666
$ws.Range("B" + $newRow).Value = 187
$ws.Range("C" + $newRow).Value = 181
$ws.Range("D" + $newRow).Value = 185
